# Applies the cryptos-list price/volume refresh described in the commit diff.
# Source cells are plain text (t="inlineStr") in the original workbook, and some
# of the replacement values look like plain numbers (e.g. "470.05"). Assigning
# such a string straight to Range.Value lets Excel auto-coerce it to a Number,
# which would silently change the cell's stored type from Text to Number.
# Set-CellText works around that by writing the value with a leading apostrophe
# (forces text entry, exactly like a user typing an apostrophe then the value in
# the Excel UI) and then resetting the cell style back to "Normal" so no
# lingering quote-prefix / text-number-format style is left attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Cell, [string]$Text)
    $range = $ws.Range($Cell)
    $looksNumeric = $Text -match '^[+-]?\d+(\.\d+)*$'
    if ($looksNumeric) {
        $range.Value = "'" + $Text
        $range.Style = "Normal"
    } else {
        $range.Value = $Text
    }
}

Set-CellText "D2" "54.798.85"
Set-CellText "E2" "  -3.13%  "
Set-CellText "D3" "2.341.77"
Set-CellText "E3" "  -5.88%  "
Set-CellText "E4" "  -0.12%  "
Set-CellText "D5" "470.05"
Set-CellText "E5" "  -4.18%  "
Set-CellText "D6" "143.58"
Set-CellText "E6" "  -2.50%  "
Set-CellText "E7" "  +0.43%  "
Set-CellText "E8" "  -2.97%  "
Set-CellText "D9" "2.338.27"
Set-CellText "E9" "  -6.79%  "
Set-CellText "D10" "0.0956"
Set-CellText "E10" "  -3.08%  "
Set-CellText "D11" "5.36"
Set-CellText "E11" "  -7.51%  "
Set-CellText "E13" "  +0.32%  "
Set-CellText "D14" "2.750.81"
Set-CellText "E14" "  -5.93%  "
Set-CellText "D15" "55.188.99"
Set-CellText "E15" "  -2.41%  "
Set-CellText "D16" "19.89"
Set-CellText "E16" "  -6.79%  "
Set-CellText "E17" "  -6.08%  "
Set-CellText "D18" "2.358.14"
Set-CellText "E18" "  -5.77%  "
Set-CellText "E19" "  -1.98%  "
Set-CellText "D20" "310.17"
Set-CellText "E20" "  -3.28%  "
Set-CellText "D21" "9.53"
Set-CellText "E21" "  -6.59%  "
Set-CellText "D22" "0.998"
Set-CellText "E22" "  +0.09%  "
Set-CellText "D23" "5.55"
Set-CellText "E23" "  -4.94%  "
Set-CellText "D24" "56.01"
Set-CellText "E24" "  -4.99%  "
Set-CellText "D25" "1.00"
Set-CellText "E25" "  +0.79%  "
Set-CellText "D26" "0.389"
Set-CellText "E26" "  -5.62%  "
Set-CellText "E27" "  -8.37%  "
Set-CellText "D28" "2.452.80"
Set-CellText "E28" "  -5.48%  "
Set-CellText "D29" "7.11"
Set-CellText "E29" "  -7.53%  "
Set-CellText "E30" "  +0.09%  "
Set-CellText "D31" "0.0₃0745"
Set-CellText "E31" "  -6.33%  "
Set-CellText "D32" "148.26"
Set-CellText "E32" "  -0.74%  "
Set-CellText "D33" "17.91"
Set-CellText "E33" "  -2.12%  "
Set-CellText "E34" "  -3.70%  "
Set-CellText "D35" "4.97"
Set-CellText "E35" "  -5.16%  "
Set-CellText "E36" "  -6.48%  "
Set-CellText "E37" "  -6.64%  "
Set-CellText "E38" "  -5.89%  "
Set-CellText "D39" "33.36"
Set-CellText "E39" "  -2.49%  "
Set-CellText "E40" "  +0.48%  "
Set-CellText "E41" "  -2.01%  "
Set-CellText "E42" "  -5.59%  "
Set-CellText "D43" "0.0940"
Set-CellText "E43" "  +1.92%  "
Set-CellText "D45" "0.570"
Set-CellText "E45" "  -7.41%  "
Set-CellText "E46" "  -0.58%  "
Set-CellText "D47" "252.03"
Set-CellText "E47" "  -2.93%  "
Set-CellText "E48" "  -4.26%  "
Set-CellText "E49" "  -9.63%  "
Set-CellText "E50" "  -5.85%  "
Set-CellText "D51" "1.758.23"
Set-CellText "E51" "  -7.70%  "
